# "Set Up XP event dispatcher (need code for addng XP still)"
#
# Changes applied:
#   1. Colour the three existing Unity bullets ("Make player attack
#      automatically.", "Split player class ... controller",
#      "Remove attack animation ...") green (RGB 00B050).
#   2. Append two new (green) runs " - Just " / " Show gizmo instead "
#      to the "Remove attack animation..." bullet.
#   3. Insert a brand-new bullet "Make xp system" right after that
#      paragraph (and before the "Notes" heading) - uncoloured.
#   4. Append a new run " - Yes" to the final "Should I create the
#      game for a player first?" bullet.

$d = $word.ActiveDocument
$green = 5287936   # RGB(0x00,0xB0,0x50) == OOXML <w:color w:val="00B050"/>

# ---------------------------------------------------------------
# 1. "Make player attack automatically."  -> green
# ---------------------------------------------------------------
$pAttack = $d.Paragraphs.Item(13)
$pAttack.Range.Font.Color = $green

# ---------------------------------------------------------------
# 2. "Split player class into player and player controller" -> green
# ---------------------------------------------------------------
$pSplit = $d.Paragraphs.Item(14)
$pSplit.Range.Font.Color = $green

# ---------------------------------------------------------------
# 3. Insert the new "Make xp system" bullet right after the
#    "Remove attack animation..." paragraph, BEFORE that paragraph
#    gets coloured/extended, so the new bullet stays uncoloured.
# ---------------------------------------------------------------
$pRemove = $d.Paragraphs.Item(15)
$pRemove.Range.InsertParagraphAfter()
$pXp = $d.Paragraphs.Item(16)
$pXp.Range.InsertAfter("Make xp system")

# ---------------------------------------------------------------
# 4. "Remove attack animation (just walking) and make a debug box
#    show instead." -> green, plus two new green runs appended:
#    " - Just " and " Show gizmo instead "
# ---------------------------------------------------------------
$pRemove = $d.Paragraphs.Item(15)
$startRemove = $pRemove.Range.Start
$endRemove = $pRemove.Range.End
$lastCharRemove = $endRemove - 1   # position right before the paragraph mark

$newRun1 = $d.Range($lastCharRemove, $lastCharRemove)
$newRun1.InsertAfter(" - Just ")

$newRun2 = $d.Range($newRun1.End, $newRun1.End)
$newRun2.InsertAfter(" Show gizmo instead ")

# Colour the original text and the two newly-typed runs.
$d.Range($startRemove, $lastCharRemove).Font.Color = $green
$d.Range($newRun1.Start, $newRun1.End).Font.Color = $green
$d.Range($newRun2.Start, $newRun2.End).Font.Color = $green

# Re-assert the colour on the paragraph's own Range: this stamps the
# paragraph-mark run properties (<w:pPr><w:rPr>) without re-merging
# the three runs we just coloured individually (they already match).
$pRemove.Range.Font.Color = $green

# ---------------------------------------------------------------
# 5. "Should I create the game for a player first?" -> append a new
#    run " - Yes"
# ---------------------------------------------------------------
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$endLast = $pLast.Range.End
$lastCharLast = $endLast - 1

$newRunYes = $d.Range($lastCharLast, $lastCharLast)
$newRunYes.InsertAfter(" - Yes")

# Toggling Bold on/off on just the new text keeps it as a run distinct
# from "Should I create the game for a player first?" instead of being
# silently coalesced back into the previous run.
$segYes = $d.Range($lastCharLast, $newRunYes.End)
$segYes.Bold = 1
$segYes.Bold = 0

Write-Host "Edit complete."
